$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Jag2"
$ws.Range("C2").Value = "Notch3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 12.97098566666667
$ws.Range("H2").Value = 38.91295700000001
$ws.Range("I2").Value = 0.7291028508134716
$ws.Range("J2").Value = 0.7291028508134717
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 3.684362666666666
$ws.Range("N2").Value = 11.053088
$ws.Range("O2").Value = 0.0397011572965827
$ws.Range("P2").Value = 0.03970115729658269
$ws.Range("Q2").Value = 47.78981534013512
$ws.Range("R2").Value = 430.108338061216
$ws.Range("S2").Value = 0.0289462269655325
$ws.Range("T2").Value = 0.0289462269655325

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Jag2"
$ws.Range("C3").Value = "Notch3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 12.97098566666667
$ws.Range("H3").Value = 38.91295700000001
$ws.Range("I3").Value = 0.7291028508134716
$ws.Range("J3").Value = 0.7291028508134717
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 5.038243666666667
$ws.Range("N3").Value = 15.114731
$ws.Range("O3").Value = 0.05429001496473517
$ws.Range("P3").Value = 0.05429001496473517
$ws.Range("Q3").Value = 65.35098638550747
$ws.Range("R3").Value = 588.1588774695672
$ws.Range("S3").Value = 0.03958300468149445
$ws.Range("T3").Value = 0.03958300468149446

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Jag2"
$ws.Range("C4").Value = "Notch3"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 12.97098566666667
$ws.Range("H4").Value = 38.91295700000001
$ws.Range("I4").Value = 0.7291028508134716
$ws.Range("J4").Value = 0.7291028508134717
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.09716666666666667
$ws.Range("N4").Value = 0.2915
$ws.Range("O4").Value = 0.001047027523164011
$ws.Range("P4").Value = 0.001047027523164011
$ws.Range("Q4").Value = 1.260347440611111
$ws.Range("R4").Value = 11.3431269655
$ws.Range("S4").Value = 0.0007633907520190488
$ws.Range("T4").Value = 0.0007633907520190489

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Jag2"
$ws.Range("C5").Value = "Notch3"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 12.97098566666667
$ws.Range("H5").Value = 38.91295700000001
$ws.Range("I5").Value = 0.7291028508134716
$ws.Range("J5").Value = 0.7291028508134717
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 83.98262666666666
$ws.Range("N5").Value = 251.94788
$ws.Range("O5").Value = 0.9049618002155182
$ws.Range("P5").Value = 0.9049618002155182
$ws.Range("Q5").Value = 1089.337446742351
$ws.Range("R5").Value = 9804.037020681162
$ws.Range("S5").Value = 0.6598102284144256
$ws.Range("T5").Value = 0.6598102284144257

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Jag2"
$ws.Range("C6").Value = "Notch3"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.047813333333333
$ws.Range("H6").Value = 3.14344
$ws.Range("I6").Value = 0.05889789011308234
$ws.Range("J6").Value = 0.05889789011308236
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 3.684362666666666
$ws.Range("N6").Value = 11.053088
$ws.Range("O6").Value = 0.0397011572965827
$ws.Range("P6").Value = 0.03970115729658269
$ws.Range("Q6").Value = 3.860524326968888
$ws.Range("R6").Value = 34.74471894272
$ws.Range("S6").Value = 0.002338314399816325
$ws.Range("T6").Value = 0.002338314399816325

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Jag2"
$ws.Range("C7").Value = "Notch3"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.047813333333333
$ws.Range("H7").Value = 3.14344
$ws.Range("I7").Value = 0.05889789011308234
$ws.Range("J7").Value = 0.05889789011308236
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 5.038243666666667
$ws.Range("N7").Value = 15.114731
$ws.Range("O7").Value = 0.05429001496473517
$ws.Range("P7").Value = 0.05429001496473517
$ws.Range("Q7").Value = 5.279138890515556
$ws.Range("R7").Value = 47.51225001464001
$ws.Range("S7").Value = 0.003197567335630568
$ws.Range("T7").Value = 0.003197567335630569

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Jag2"
$ws.Range("C8").Value = "Notch3"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.047813333333333
$ws.Range("H8").Value = 3.14344
$ws.Range("I8").Value = 0.05889789011308234
$ws.Range("J8").Value = 0.05889789011308236
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.09716666666666667
$ws.Range("N8").Value = 0.2915
$ws.Range("O8").Value = 0.001047027523164011
$ws.Range("P8").Value = 0.001047027523164011
$ws.Range("Q8").Value = 0.1018125288888889
$ws.Range("R8").Value = 0.91631276
$ws.Range("S8").Value = 0.00006166771200468671
$ws.Range("T8").Value = 0.00006166771200468673

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Jag2"
$ws.Range("C9").Value = "Notch3"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.047813333333333
$ws.Range("H9").Value = 3.14344
$ws.Range("I9").Value = 0.05889789011308234
$ws.Range("J9").Value = 0.05889789011308236
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 83.98262666666666
$ws.Range("N9").Value = 251.94788
$ws.Range("O9").Value = 0.9049618002155182
$ws.Range("P9").Value = 0.9049618002155182
$ws.Range("Q9").Value = 87.99811598968888
$ws.Range("R9").Value = 791.9830439072
$ws.Range("S9").Value = 0.05330034066563077
$ws.Range("T9").Value = 0.05330034066563078

# Row 10
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Jag2"
$ws.Range("C10").Value = "Notch3"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.312552
$ws.Range("H10").Value = 3.937656
$ws.Range("I10").Value = 0.07377892703252469
$ws.Range("J10").Value = 0.0737789270325247
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 3.684362666666666
$ws.Range("N10").Value = 11.053088
$ws.Range("O10").Value = 0.0397011572965827
$ws.Range("P10").Value = 0.03970115729658269
$ws.Range("Q10").Value = 4.835917586858666
$ws.Range("R10").Value = 43.523258281728
$ws.Range("S10").Value = 0.00292910878729136
$ws.Range("T10").Value = 0.00292910878729136

# Row 11
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Jag2"
$ws.Range("C11").Value = "Notch3"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 1.312552
$ws.Range("H11").Value = 3.937656
$ws.Range("I11").Value = 0.07377892703252469
$ws.Range("J11").Value = 0.0737789270325247
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 5.038243666666667
$ws.Range("N11").Value = 15.114731
$ws.Range("O11").Value = 0.05429001496473517
$ws.Range("P11").Value = 0.05429001496473517
$ws.Range("Q11").Value = 6.612956801170667
$ws.Range("R11").Value = 59.51661121053601
$ws.Range("S11").Value = 0.004005459052677869
$ws.Range("T11").Value = 0.00400545905267787

# Row 12
$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Jag2"
$ws.Range("C12").Value = "Notch3"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 1.312552
$ws.Range("H12").Value = 3.937656
$ws.Range("I12").Value = 0.07377892703252469
$ws.Range("J12").Value = 0.0737789270325247
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.09716666666666667
$ws.Range("N12").Value = 0.2915
$ws.Range("O12").Value = 0.001047027523164011
$ws.Range("P12").Value = 0.001047027523164011
$ws.Range("Q12").Value = 0.1275363026666667
$ws.Range("R12").Value = 1.147826724
$ws.Range("S12").Value = 0.00007724856723256264
$ws.Range("T12").Value = 0.00007724856723256265

# Row 13
$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Jag2"
$ws.Range("C13").Value = "Notch3"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 1.312552
$ws.Range("H13").Value = 3.937656
$ws.Range("I13").Value = 0.07377892703252469
$ws.Range("J13").Value = 0.0737789270325247
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 83.98262666666666
$ws.Range("N13").Value = 251.94788
$ws.Range("O13").Value = 0.9049618002155182
$ws.Range("P13").Value = 0.9049618002155182
$ws.Range("Q13").Value = 110.2315645965867
$ws.Range("R13").Value = 992.08408136928
$ws.Range("S13").Value = 0.0667671106253229
$ws.Range("T13").Value = 0.06676711062532291

# Row 14
$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Jag2"
$ws.Range("C14").Value = "Notch3"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 2.458986333333333
$ws.Range("H14").Value = 7.376958999999999
$ws.Range("I14").Value = 0.1382203320409214
$ws.Range("J14").Value = 0.1382203320409214
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 3.684362666666666
$ws.Range("N14").Value = 11.053088
$ws.Range("O14").Value = 0.0397011572965827
$ws.Range("P14").Value = 0.03970115729658269
$ws.Range("Q14").Value = 9.059797444376889
$ws.Range("R14").Value = 81.53817699939198
$ws.Range("S14").Value = 0.005487507143942508
$ws.Range("T14").Value = 0.005487507143942508

# Row 15
$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Jag2"
$ws.Range("C15").Value = "Notch3"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 2.458986333333333
$ws.Range("H15").Value = 7.376958999999999
$ws.Range("I15").Value = 0.1382203320409214
$ws.Range("J15").Value = 0.1382203320409214
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 5.038243666666667
$ws.Range("N15").Value = 15.114731
$ws.Range("O15").Value = 0.05429001496473517
$ws.Range("P15").Value = 0.05429001496473517
$ws.Range("Q15").Value = 12.38897232033656
$ws.Range("R15").Value = 111.500750883029
$ws.Range("S15").Value = 0.007503983894932284
$ws.Range("T15").Value = 0.007503983894932286

# Row 16
$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Jag2"
$ws.Range("C16").Value = "Notch3"
$ws.Range("D16").Value = "M2"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 2.458986333333333
$ws.Range("H16").Value = 7.376958999999999
$ws.Range("I16").Value = 0.1382203320409214
$ws.Range("J16").Value = 0.1382203320409214
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.09716666666666667
$ws.Range("N16").Value = 0.2915
$ws.Range("O16").Value = 0.001047027523164011
$ws.Range("P16").Value = 0.001047027523164011
$ws.Range("Q16").Value = 0.2389315053888889
$ws.Range("R16").Value = 2.1503835485
$ws.Range("S16").Value = 0.0001447204919077131
$ws.Range("T16").Value = 0.0001447204919077131

# Row 17
$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "Jag2"
$ws.Range("C17").Value = "Notch3"
$ws.Range("D17").Value = "sCs"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 2.458986333333333
$ws.Range("H17").Value = 7.376958999999999
$ws.Range("I17").Value = 0.1382203320409214
$ws.Range("J17").Value = 0.1382203320409214
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 83.98262666666666
$ws.Range("N17").Value = 251.94788
$ws.Range("O17").Value = 0.9049618002155182
$ws.Range("P17").Value = 0.9049618002155182
$ws.Range("Q17").Value = 206.5121312107689
$ws.Range("R17").Value = 1858.60918089692
$ws.Range("S17").Value = 0.1250841205101389
$ws.Range("T17").Value = 0.1250841205101389
